$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = 0.7825479339666588
$ws.Range("M2").Value = 13.67700833333333
$ws.Range("N2").Value = 41.031025
$ws.Range("O2").Value = 0.124413831206147
$ws.Range("P2").Value = 0.124413831206147
$ws.Range("Q2").Value = 1.539360964925
$ws.Range("R2").Value = 13.854248684325
$ws.Range("S2").Value = 0.09735978656724693
$ws.Range("T2").Value = 0.09735978656724692
$ws.Range("J3").Value = 0.7825479339666588
$ws.Range("M3").Value = 74.64939600000001
$ws.Range("O3").Value = 0.679053278848249
$ws.Range("P3").Value = 0.6790532788482488
$ws.Range("R3").Value = 75.616777522764
$ws.Range("S3").Value = 0.5313917404159827
$ws.Range("T3").Value = 0.5313917404159826
$ws.Range("J4").Value = 0.7825479339666588
$ws.Range("M4").Value = 1.629335666666667
$ws.Range("N4").Value = 4.888007
$ws.Range("O4").Value = 0.01482136207497777
$ws.Range("P4").Value = 0.01482136207497777
$ws.Range("Q4").Value = 0.183383358619
$ws.Range("R4").Value = 1.650450227571
$ws.Range("S4").Value = 0.01159842627034565
$ws.Range("T4").Value = 0.01159842627034565
$ws.Range("J5").Value = 0.7825479339666588
$ws.Range("M5").Value = 19.17462033333333
$ws.Range("N5").Value = 57.523861
$ws.Range("O5").Value = 0.174423230537864
$ws.Range("P5").Value = 0.174423230537864
$ws.Range("Q5").Value = 2.158122693137
$ws.Range("R5").Value = 19.423104238233
$ws.Range("S5").Value = 0.1364945386931957
$ws.Range("T5").Value = 0.1364945386931957
$ws.Range("J6").Value = 0.7825479339666588
$ws.Range("M6").Value = 0.801214
$ws.Range("N6").Value = 2.403642
$ws.Range("O6").Value = 0.007288297332762355
$ws.Range("P6").Value = 0.007288297332762355
$ws.Range("Q6").Value = 0.090177436914
$ws.Range("R6").Value = 0.8115969322259999
$ws.Range("S6").Value = 0.005703442019887891
$ws.Range("T6").Value = 0.00570344201988789
$ws.Range("G7").Value = 0.03127533333333334
$ws.Range("H7").Value = 0.09382600000000001
$ws.Range("I7").Value = 0.2174520660333412
$ws.Range("J7").Value = 0.2174520660333412
$ws.Range("M7").Value = 13.67700833333333
$ws.Range("N7").Value = 41.031025
$ws.Range("O7").Value = 0.124413831206147
$ws.Range("P7").Value = 0.124413831206147
$ws.Range("Q7").Value = 0.4277529946277778
$ws.Range("R7").Value = 3.84977695165
$ws.Range("S7").Value = 0.02705404463890003
$ws.Range("T7").Value = 0.02705404463890003
$ws.Range("G8").Value = 0.03127533333333334
$ws.Range("H8").Value = 0.09382600000000001
$ws.Range("I8").Value = 0.2174520660333412
$ws.Range("J8").Value = 0.2174520660333412
$ws.Range("M8").Value = 74.64939600000001
$ws.Range("O8").Value = 0.679053278848249
$ws.Range("P8").Value = 0.6790532788482488
$ws.Range("Q8").Value = 2.334684743032001
$ws.Range("S8").Value = 0.1476615384322663
$ws.Range("T8").Value = 0.1476615384322662
$ws.Range("G9").Value = 0.03127533333333334
$ws.Range("H9").Value = 0.09382600000000001
$ws.Range("I9").Value = 0.2174520660333412
$ws.Range("J9").Value = 0.2174520660333412
$ws.Range("M9").Value = 1.629335666666667
$ws.Range("N9").Value = 4.888007
$ws.Range("O9").Value = 0.01482136207497777
$ws.Range("P9").Value = 0.01482136207497777
$ws.Range("Q9").Value = 0.05095801608688889
$ws.Range("R9").Value = 0.458622144782
$ws.Range("S9").Value = 0.003222935804632125
$ws.Range("T9").Value = 0.003222935804632125
$ws.Range("G10").Value = 0.03127533333333334
$ws.Range("H10").Value = 0.09382600000000001
$ws.Range("I10").Value = 0.2174520660333412
$ws.Range("J10").Value = 0.2174520660333412
$ws.Range("M10").Value = 19.17462033333333
$ws.Range("N10").Value = 57.523861
$ws.Range("O10").Value = 0.174423230537864
$ws.Range("P10").Value = 0.174423230537864
$ws.Range("Q10").Value = 0.5996926424651111
$ws.Range("R10").Value = 5.397233782186
$ws.Range("S10").Value = 0.03792869184466829
$ws.Range("T10").Value = 0.03792869184466829
$ws.Range("G11").Value = 0.03127533333333334
$ws.Range("H11").Value = 0.09382600000000001
$ws.Range("I11").Value = 0.2174520660333412
$ws.Range("J11").Value = 0.2174520660333412
$ws.Range("M11").Value = 0.801214
$ws.Range("N11").Value = 2.403642
$ws.Range("O11").Value = 0.007288297332762355
$ws.Range("P11").Value = 0.007288297332762355
$ws.Range("Q11").Value = 0.02505823492133333
$ws.Range("R11").Value = 0.225524114292
$ws.Range("S11").Value = 0.001584855312874464
$ws.Range("T11").Value = 0.001584855312874464
